$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of IRS EIN-assignment data appended to the export (rows 183-190).
# Force the newly-written range to Text format first so EIN/SSN/ZIP values
# that look numeric (and especially ones with leading zeros, e.g. "099664210")
# are preserved exactly as text, matching the rest of the sheet.
$newRange = $ws.Range("A183:J190")
$newRange.NumberFormat = "@"

# Row 183
$ws.Range("A183").Value = '33-4000359'
$ws.Range("B183").Value = 'BLUESKY INNOVATIONS LLC'
$ws.Range("C183").Value = '432813777'
$ws.Range("D183").Value = 'Basilia'
$ws.Range("E183").Value = 'Gonzalez'
$ws.Range("F183").Value = '1992-01-02 00:00:00'
$ws.Range("G183").Value = '1828 James Street'
$ws.Range("H183").Value = 'Irving'
$ws.Range("I183").Value = 'TX'
$ws.Range("J183").Value = '75061'

# Row 184
$ws.Range("A184").Value = '33-4000640'
$ws.Range("B184").Value = 'STELLARPEAK VENTURES LLC'
$ws.Range("C184").Value = '099664210'
$ws.Range("D184").Value = 'RAFAEL'
$ws.Range("E184").Value = 'PEREZ-ESPEJO'
$ws.Range("F184").Value = '1958-10-27 00:00:00'
$ws.Range("G184").Value = '1100 GOUGH ST'
$ws.Range("H184").Value = 'SAN FRANCISCO'
$ws.Range("I184").Value = 'CA'
$ws.Range("J184").Value = '94109'

# Row 185
$ws.Range("A185").Value = '33-4000781'
$ws.Range("B185").Value = 'NEXUSFLOW TECHNOLOGIES LLC'
$ws.Range("C185").Value = '643015000'
$ws.Range("D185").Value = 'ZACHARY'
$ws.Range("E185").Value = 'DAVIDSON'
$ws.Range("F185").Value = '1980-09-17 00:00:00'
$ws.Range("G185").Value = '955 PINE STREET'
$ws.Range("H185").Value = 'SAN FRANCISCO'
$ws.Range("I185").Value = 'CA'
$ws.Range("J185").Value = '94108'

# Row 186
$ws.Range("A186").Value = '33-4000836'
$ws.Range("B186").Value = 'APEXVIBE GLOBAL LLC'
$ws.Range("C186").Value = '625247624'
$ws.Range("D186").Value = 'JOSHUA'
$ws.Range("E186").Value = 'BONGAWIL'
$ws.Range("F186").Value = '1988-11-02 00:00:00'
$ws.Range("G186").Value = '236 SAGEBRUSH LANE'
$ws.Range("H186").Value = 'AMERICAN CANYON'
$ws.Range("I186").Value = 'CA'
$ws.Range("J186").Value = '94503'

# Row 187
$ws.Range("A187").Value = '33-4000904'
$ws.Range("B187").Value = 'HORIZONWAVE INDUSTRIES LLC'
$ws.Range("C187").Value = '620425729'
$ws.Range("D187").Value = 'ARACELI'
$ws.Range("E187").Value = 'MARTIN'
$ws.Range("F187").Value = '1990-01-09 00:00:00'
$ws.Range("G187").Value = '190 HALE ST'
$ws.Range("H187").Value = 'SAN FRANCISCO'
$ws.Range("I187").Value = 'CA'
$ws.Range("J187").Value = '94134'

# Row 188
$ws.Range("A188").Value = '33-4000992'
$ws.Range("B188").Value = 'REDSTONE STRATEGIES LLC'
$ws.Range("C188").Value = '615361434'
$ws.Range("D188").Value = 'HARMONIE'
$ws.Range("E188").Value = 'WONG'
$ws.Range("F188").Value = '1990-03-04 00:00:00'
$ws.Range("G188").Value = '42 BERKELEY WAY'
$ws.Range("H188").Value = 'SAN FRANCISCO'
$ws.Range("I188").Value = 'CA'
$ws.Range("J188").Value = '94131'

# Row 189
$ws.Range("A189").Value = '33-4001063'
$ws.Range("B189").Value = 'NOVACRAFT HOLDINGS LLC'
$ws.Range("C189").Value = '612407257'
$ws.Range("D189").Value = 'JENNY'
$ws.Range("E189").Value = 'VALDEZ'
$ws.Range("F189").Value = '1988-03-07 00:00:00'
$ws.Range("G189").Value = '990 MAGNOLIA AVE'
$ws.Range("H189").Value = 'MILLBRAE'
$ws.Range("I189").Value = 'CA'
$ws.Range("J189").Value = '94030'

# Row 190
$ws.Range("A190").Value = '33-4001534'
$ws.Range("B190").Value = 'CRYSTALCOVE ENTERPRISES LLC'
$ws.Range("C190").Value = '630045084'
$ws.Range("D190").Value = 'AVIYANCA'
$ws.Range("E190").Value = 'PRAKASH'
$ws.Range("F190").Value = '1993-11-19 00:00:00'
$ws.Range("G190").Value = '711 TORREYA AVENUE'
$ws.Range("H190").Value = 'SUNNYVALE'
$ws.Range("I190").Value = 'CA'
$ws.Range("J190").Value = '94086'
